# Updated for seperate the P2P process
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BatchSheet")

# Rename the "Cloud.*" actions to "P2P.*" for the rows that moved into the
# dedicated P2P process group (rows 2, 3, 8, 9, 10 -> column E / "Action4").
$ws.Range("E2").Value  = "P2P.createItemBasedRequisition"
$ws.Range("E3").Value  = "P2P.createPurchaseOrder"
$ws.Range("E8").Value  = "P2P.createReceivingReceipt"
$ws.Range("E9").Value  = "P2P.createPurOrderMatchedInvoice"
$ws.Range("E10").Value = "P2P.createPaymentQuickCheck"

# Update the saved selection/active cell on the sheet to match the new state.
$ws.Activate()
$ws.Range("E17").Select()
